$d = $word.ActiveDocument

# A manual line break (OOXML <w:br/>) shows up in the Word Range/Find text
# stream as Chr(11) ("manual line break" character).
$lb = [char]11

$movedText = "LOM3231 -  Métodos Experimentais da Física IV  (Indicação de Conjunto)" + $lb

# Locate the run (text + trailing line break) that needs to move: the
# "LOM3231 ... (Indicação de Conjunto)" requirement entry.
$moveRange = $d.Content
$foundMove = $moveRange.Find.Execute($movedText, $false, $false, $false, $false, $false, `
                                      $true, 1, $false, "", 0)
if (-not $foundMove) {
    throw "Could not find the LOM3231 requirement entry to move"
}
$moveStart = $moveRange.Start
$moveEnd = $moveRange.End

# Locate the insertion point: right before the "LOM3206" entry, which is
# where the LOM3231 entry needs to end up (as the new first item).
$targetRange = $d.Content
$foundTarget = $targetRange.Find.Execute("LOM3206 -  Eletrônica  (Requisito)", $false, $false, `
                                          $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundTarget) {
    throw "Could not find the LOM3206 requirement entry to insert before"
}
$targetStart = $targetRange.Start

# Remove the LOM3231 entry (text + w:br) from its current location, after
# the LOM3215 entry.
$deleteRange = $d.Range($moveStart, $moveEnd)
$deleteRange.Delete()

# Re-insert the same text (with its line break) immediately before the
# LOM3206 entry. $targetStart is still valid: it is upstream of the range
# that was just deleted, so nothing shifted its offset.
$insertRange = $d.Range($targetStart, $targetStart)
$insertRange.InsertBefore($movedText)

Write-Output "Moved LOM3231 requirement entry before LOM3206"
